# T02-Solenoid_driver_unit_test.docx — "Updated Test Plan document"
#
# The underlying edit was: the author's cursor/last-edit position (tracked
# by Word's hidden "_GoBack" bookmark) moved from the very start of the
# document (just before the title text) down into the test-procedure
# table, landing between "supply," and " after approximately 10ms." in
# step 4's expected-result cell. Word only ever keeps a single "_GoBack"
# bookmark, so the old one is removed and a new one is dropped at the new
# location, which splits that run of text in two.

$d = $word.ActiveDocument

# 1) Remove the old "_GoBack" bookmark sitting at the top of the document
#    (inside the title paragraph's properties, right before the title run).
$goBack = $d.Bookmarks.Item("_GoBack")
if ($goBack -ne $null) {
    $goBack.Delete()
}

# 2) Locate "supply, after approximately 10ms." in the procedure table and
#    drop a fresh "_GoBack" bookmark right after "supply," (i.e. before the
#    space that starts " after approximately 10ms."). Adding a bookmark at
#    that collapsed range splits the run in two, exactly like Word does.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("supply, after approximately 10ms.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $splitPoint = $searchRange.Start + 7   # length of "supply," (7 chars)
    $insertionRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $insertionRange)
}
